$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column (D) cells being updated, to prevent Excel
# from auto-converting numeric-looking strings into numbers (losing trailing
# zeros / exact text representation), matching the original inline-string text cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.044.45"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.794.12"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "699.38"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.81"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.794.41"
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.52"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.15"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.435.98"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.813.11"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.130.69"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.54"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "512.38"
$ws.Range("E21").Value = "  +2.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.40"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.713"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.46"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000140"
$ws.Range("E25").Value = "  -4.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.60"
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.942.85"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.98"
$ws.Range("E30").Value = "  -5.76%  "
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.26"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.27"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.12"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.31"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.757.47"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.67"
$ws.Range("E39").Value = "  +10.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.35"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.16"
$ws.Range("E44").Value = "  -6.02%  "
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.27"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.26"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000300"
$ws.Range("E48").Value = "  -6.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "422.59"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.38"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.63"
$ws.Range("E51").Value = "  -1.66%  "
